# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    3  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    4  = @(0.127881588408715, 0.3127903958511391, 26.21740644021617, 8.660232485948974, 35.318310910425)
    5  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    6  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    7  = @(0.6753301551942219, 10.29869402782916, 337.1190423067083, 645.3272768299601, 993.4203433196917)
    8  = @(0.04763786555579896, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0.7443468554461139)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    10 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    11 = @(0.3048080303191223, 1.667794583268128, 3.900430680208489, 8.660232485948974, 14.53326577974471)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
